$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 3127.238
$ws.Range("I135").Value = 3135.4211
$ws.Range("K135").Value = 28218.7899
$ws.Range("M135").Value = -25683.7899
$ws.Range("H138").Value = 2893.1
$ws.Range("J138").Value = 3187.014
$ws.Range("L138").Value = 9561.042000000001
$ws.Range("N138").Value = -19841.042

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1115.2667
$ws.Range("J2").Value = 912.5
$ws.Range("L2").Value = 912.5
$ws.Range("N2").Value = -1138.5
$ws.Range("H28").Value = 46630
$ws.Range("I28").Value = 14945
$ws.Range("K28").Value = 14945
$ws.Range("M28").Value = -14753
$ws.Range("H61").Value = 14744049
$ws.Range("I61").Value = 21743846
$ws.Range("J61").Value = 108112.18
$ws.Range("K61").Value = 21743846
$ws.Range("L61").Value = 108112.18
$ws.Range("M61").Value = -21743634
$ws.Range("N61").Value = -108536.18
$ws.Range("H74").Value = 11373918
$ws.Range("I74").Value = 14709253
$ws.Range("J74").Value = 33780.4
$ws.Range("K74").Value = 14709253
$ws.Range("L74").Value = 33780.4
$ws.Range("M74").Value = -14708379
$ws.Range("N74").Value = -35528.4
$ws.Range("H77").Value = 11373918
$ws.Range("I77").Value = 14709253
$ws.Range("J77").Value = 33780.4
$ws.Range("K77").Value = 73546265
$ws.Range("L77").Value = 168902
$ws.Range("M77").Value = -73541897
$ws.Range("N77").Value = -177638
$ws.Range("H99").Value = 46630
$ws.Range("I99").Value = 14945
$ws.Range("K99").Value = 14945
$ws.Range("M99").Value = -11950
$ws.Range("H116").Value = 1115.2667
$ws.Range("J116").Value = 912.5
$ws.Range("L116").Value = 912.5
$ws.Range("N116").Value = -5500.5
$ws.Range("H136").Value = 14744049
$ws.Range("I136").Value = 21743846
$ws.Range("J136").Value = 108112.18
$ws.Range("K136").Value = 65231538
$ws.Range("L136").Value = 324336.54
$ws.Range("M136").Value = -65228988
$ws.Range("N136").Value = -329436.54

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1115.2667
$ws.Range("J3").Value = 912.5
$ws.Range("L3").Value = 912.5
$ws.Range("N3").Value = -1140.5
$ws.Range("H98").Value = 76316.2
$ws.Range("J98").Value = 76316.2
$ws.Range("L98").Value = 76316.2
$ws.Range("N98").Value = -82306.2
$ws.Range("H105").Value = 2362.111
$ws.Range("I105").Value = 2177.5715
$ws.Range("K105").Value = 2177.5715
$ws.Range("M105").Value = -430.5715
$ws.Range("H133").Value = 47399.8
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 3000
$ws.Range("J25").Value = 3000
$ws.Range("L25").Value = 3000
$ws.Range("N25").Value = -3348
$ws.Range("H41").Value = 10000
$ws.Range("I41").Value = 10000
$ws.Range("K41").Value = 10000
$ws.Range("M41").Value = -9572
$ws.Range("H50").Value = 39999
$ws.Range("I50").Value = 39999
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 39999
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -39374
$ws.Range("N50").ClearContents()
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = 0
$ws.Range("H60").Value = 56697.668
$ws.Range("I60").Value = 30093
$ws.Range("J60").Value = 70000
$ws.Range("K60").Value = 30093
$ws.Range("L60").Value = 70000
$ws.Range("M60").Value = -29582
$ws.Range("N60").Value = -71022
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = 0
$ws.Range("H132").Value = 2043.7709
$ws.Range("I132").Value = 1830.3489
$ws.Range("K132").Value = 5491.0467
$ws.Range("M132").Value = -2961.0467
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3615.3635
$ws.Range("I80").Value = 3347.5
$ws.Range("J80").Value = 3674.889
$ws.Range("K80").Value = 10042.5
$ws.Range("L80").Value = 11024.667
$ws.Range("M80").Value = -9106.5
$ws.Range("N80").Value = -12896.667
$ws.Range("H83").Value = 3615.3635
$ws.Range("I83").Value = 3347.5
$ws.Range("J83").Value = 3674.889
$ws.Range("K83").Value = 30127.5
$ws.Range("L83").Value = 33074.001
$ws.Range("M83").Value = -25447.5
$ws.Range("N83").Value = -42434.001
$ws.Range("H92").Value = 627649.3
$ws.Range("J92").Value = 2970.7144
$ws.Range("L92").Value = 8912.143199999999
$ws.Range("N92").Value = -11408.1432
$ws.Range("H131").Value = 11123.105
$ws.Range("I131").Value = 13124.75
$ws.Range("J131").Value = 10589.333
$ws.Range("K131").Value = 39374.25
$ws.Range("L131").Value = 31767.999
$ws.Range("M131").Value = -34334.25
$ws.Range("N131").Value = -41847.999
$ws.Range("H132").Value = 2119.3635
$ws.Range("I132").Value = 1996.3077
$ws.Range("J132").Value = 2297.111
$ws.Range("K132").Value = 17966.7693
$ws.Range("L132").Value = 20673.999
$ws.Range("M132").Value = -15436.7693
$ws.Range("N132").Value = -25733.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 89988.25
$ws.Range("I5").Value = 89988.25
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 89988.25
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -89876.25
$ws.Range("N5").ClearContents()
$ws.Range("H132").Value = 33338942
$ws.Range("I132").Value = 43479770
$ws.Range("J132").Value = 19080.428
$ws.Range("K132").Value = 130439310
$ws.Range("L132").Value = 57241.284
$ws.Range("M132").Value = -130436780
$ws.Range("N132").Value = -62301.284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1656.6666
$ws.Range("I61").Value = 1844.5714
$ws.Range("J61").Value = 999
$ws.Range("K61").Value = 1844.5714
$ws.Range("L61").Value = 999
$ws.Range("M61").Value = -1642.5714
$ws.Range("N61").Value = -1403
$ws.Range("H81").Value = 299
$ws.Range("I81").Value = 299
$ws.Range("K81").Value = 598
$ws.Range("M81").Value = 463
$ws.Range("H84").Value = 299
$ws.Range("I84").Value = 299
$ws.Range("K84").Value = 2990
$ws.Range("M84").Value = 2314
$ws.Range("H99").Value = 61000
$ws.Range("I99").Value = 61000
$ws.Range("K99").Value = 61000
$ws.Range("M99").Value = -58005
$ws.Range("H113").Value = 1656.6666
$ws.Range("I113").Value = 1844.5714
$ws.Range("J113").Value = 999
$ws.Range("K113").Value = 1844.5714
$ws.Range("L113").Value = 999
$ws.Range("M113").Value = 325.4286
$ws.Range("N113").Value = -5339
$ws.Range("H136").Value = 49167.594
$ws.Range("I136").Value = 7268.5293
$ws.Range("J136").Value = 120396
$ws.Range("K136").Value = 21805.5879
$ws.Range("L136").Value = 361188
$ws.Range("M136").Value = -19255.5879
$ws.Range("N136").Value = -366288

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1482.9788
$ws.Range("I132").Value = 1509.1282
$ws.Range("K132").Value = 4527.3846
$ws.Range("M132").Value = -1997.3846
$ws.Range("H136").Value = 3159.8
$ws.Range("I136").Value = 1300
$ws.Range("K136").Value = 3900
$ws.Range("M136").Value = -1350
